$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... на 13 вопросов." -> "... на 14 вопросов."
#    The "1" and "3" of "13" live in two separate runs with a collapsed
#    "_GoBack" bookmark sitting right after the "3". Touch only the single
#    character "3" so the bookmark and the neighbouring "1" run survive.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx13 = $full.IndexOf("на 13 вопросов")
$idx3 = $idx13 + 4
$digit = $d.Range($idx3, $idx3 + 1)
if ($digit.Text -eq "3") {
    $digit.Text = "4"
}

# ---------------------------------------------------------------------------
# 2) Expand the last Q&A block:
#      "Результат: ... данную идею. "   (unchanged text, now one flowing run)
#      <blank paragraph>
#      "Поддержали бы вы ... пользователь?"   (new, numbered + italic question)
#      <blank paragraph>
#      "Результат: Практически ... личностей. "   (new answer paragraph)
#      <blank paragraph>
# ---------------------------------------------------------------------------

# Locate the paragraph that currently reads
# "Результат: Все респонденты ответили, что они полностью поддерживают данную идею. "
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt.StartsWith("Результат") -and $ptxt.Contains("полностью поддерживают")) {
        $target = $d.Paragraphs($i)
        break
    }
}

$targetIndex = $target.Index

# Re-flow the tail of the paragraph (everything after the underlined
# "Результат" label) into a single run. The diff merges the two runs
# ": Все респонденты ответили, что " + "они полностью поддерживают данную
# идею" into one (the trailing ". " run already matches and folds in too),
# while leaving the underlined "Результат" run untouched.
$target.Range.Find.Execute(
    ": Все респонденты ответили, что они полностью поддерживают данную идею",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ": Все респонденты ответили, что они полностью поддерживают данную идею", 2)

# Insert five new paragraphs right after it (all starting out as plain,
# non-numbered / non-italic clones of $target's paragraph formatting).
for ($i = 0; $i -lt 5; $i++) {
    $target.Range.InsertParagraphAfter()
}

# New paragraph #2 after target -> the numbered, italic question.
$qIndex = $targetIndex + 2
$q = $d.Paragraphs($qIndex)
$q.Range.Text = "Поддержали бы вы других пользователей материальными средствами для продвижения контента, который создает данный пользователь?"
$q.Range.Font.Italic = 1

$listSource = $d.Paragraphs($qIndex - 4)
$q.Range.ListFormat.ListTemplate = $listSource.Range.ListFormat.ListTemplate
$q.Range.ListFormat.ListLevelNumber = $listSource.Range.ListFormat.ListLevelNumber

# New paragraph #4 after target -> the new "Результат" answer.
$aIndex = $targetIndex + 4
$a = $d.Paragraphs($aIndex)
$a.Range.Text = "Результат: Практически все согласились с данным предложением и с тем, что необходимо поощрять творческих личностей. "
$aStart = $a.Range.Start
$aURange = $d.Range($aStart, $aStart + 9)
$aURange.Font.Underline = 1
